{"js": "// Replace the three-digit-by-one-digit multiplication equations whose\n// results changed in this revision. Each original equation string is\n// unique in the document, so a simple exact-text search/replace per pair\n// is safe and unambiguous.\nconst replacements = [\n  [\"346\u00d78=2768\", \"889\u00d77=6223\"],\n  [\"262\u00d73=786\", \"172\u00d79=1548\"],\n  [\"739\u00d78=5912\", \"723\u00d78=5784\"],\n  [\"576\u00d75=2880\", \"540\u00d78=4320\"],\n  [\"702\u00d79=6318\", \"376\u00d77=2632\"],\n  [\"817\u00d77=5719\", \"778\u00d73=2334\"],\n  [\"417\u00d77=2919\", \"166\u00d72=332\"],\n  [\"584\u00d72=1168\", \"630\u00d72=1260\"],\n  [\"513\u00d75=2565\", \"831\u00d77=5817\"],\n  [\"580\u00d74=2320\", \"422\u00d79=3798\"],\n  [\"978\u00d77=6846\", \"963\u00d79=8667\"],\n  [\"742\u00d72=1484\", \"659\u00d77=4613\"],\n  [\"437\u00d75=2185\", \"495\u00d72=990\"],\n  [\"961\u00d72=1922\", \"432\u00d76=2592\"],\n  [\"151\u00d77=1057\", \"279\u00d77=1953\"],\n  [\"605\u00d75=3025\", \"609\u00d72=1218\"],\n  [\"212\u00d74=848\", \"113\u00d74=452\"],\n  [\"333\u00d73=999\", \"201\u00d72=402\"],\n  [\"182\u00d74=728\", \"396\u00d73=1188\"],\n  [\"215\u00d79=1935\", \"325\u00d79=2925\"],\n  [\"810\u00d79=7290\", \"434\u00d72=868\"],\n  [\"908\u00d76=5448\", \"510\u00d75=2550\"],\n  [\"684\u00d72=1368\", \"227\u00d77=1589\"],\n  [\"458\u00d79=4122\", \"928\u00d77=6496\"],\n  [\"909\u00d73=2727\", \"792\u00d79=7128\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit-by-one-digit multiplication equations in the\n# results table. Every \"old\" equation string is unique in the document, so\n# a Find/Replace-All per pair unambiguously lands on the correct cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"346\u00d78=2768\", \"889\u00d77=6223\"),\n    @(\"262\u00d73=786\",  \"172\u00d79=1548\"),\n    @(\"739\u00d78=5912\", \"723\u00d78=5784\"),\n    @(\"576\u00d75=2880\", \"540\u00d78=4320\"),\n    @(\"702\u00d79=6318\", \"376\u00d77=2632\"),\n    @(\"817\u00d77=5719\", \"778\u00d73=2334\"),\n    @(\"417\u00d77=2919\", \"166\u00d72=332\"),\n    @(\"584\u00d72=1168\", \"630\u00d72=1260\"),\n    @(\"513\u00d75=2565\", \"831\u00d77=5817\"),\n    @(\"580\u00d74=2320\", \"422\u00d79=3798\"),\n    @(\"978\u00d77=6846\", \"963\u00d79=8667\"),\n    @(\"742\u00d72=1484\", \"659\u00d77=4613\"),\n    @(\"437\u00d75=2185\", \"495\u00d72=990\"),\n    @(\"961\u00d72=1922\", \"432\u00d76=2592\"),\n    @(\"151\u00d77=1057\", \"279\u00d77=1953\"),\n    @(\"605\u00d75=3025\", \"609\u00d72=1218\"),\n    @(\"212\u00d74=848\",  \"113\u00d74=452\"),\n    @(\"333\u00d73=999\",  \"201\u00d72=402\"),\n    @(\"182\u00d74=728\",  \"396\u00d73=1188\"),\n    @(\"215\u00d79=1935\", \"325\u00d79=2925\"),\n    @(\"810\u00d79=7290\", \"434\u00d72=868\"),\n    @(\"908\u00d76=5448\", \"510\u00d75=2550\"),\n    @(\"684\u00d72=1368\", \"227\u00d77=1589\"),\n    @(\"458\u00d79=4122\", \"928\u00d77=6496\"),\n    @(\"909\u00d73=2727\", \"792\u00d79=7128\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
